# BB_Build.xlsx - SLG building config: add Icon / ShowName fields.
# Two new rows are inserted ahead of the existing "Desc" row (which shifts
# from row 7 down to row 9), then populated with the new field metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 7 and 8, pushing the old row 7 ("Desc") to row 9.
$ws.Rows("7:8").Insert()

# Row 7: new "Icon" field definition.
$ws.Range("A7").Value = "Icon"
$ws.Range("B7").Value = "string"
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = $false
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = "Friend"
$ws.Range("J7").Value = "图标"

# Row 8: new "ShowName" field definition.
$ws.Range("A8").Value = "ShowName"
$ws.Range("B8").Value = "string"
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = $false
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = "Friend"
$ws.Range("J8").Value = "名字"

# Match the author's recorded selection at save time.
[void]$ws.Range("C13").Select()
